# Update the "Case locations and outbreaks - case alerts / public exposure sites" worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old "Hoppers Crossing / Caltex Woolworths" row (previously row 15).
#    Everything below it shifts up by one row.
$ws.Rows.Item(15).Delete()

# 2. Insert four new blank rows after the "Maidstone" row (now row 17), before
#    "Melbourne Airport" (now row 18), to hold the new Melbourne venue entries.
$ws.Rows.Item(18).Resize(4).Insert()

# 3. Populate the four newly inserted rows (18-21) with the new exposure sites.
$newRows = @(
    @(18, "Melbourne", "Queen Victoria Market  Queen Street  Melbourne VIC 3000", "8:25am - 10:10am  11/2/2021", "Case attended Section 2 - Fruit and Vegetables, and used Section 2 female toilets. See a map of the Queen Victoria Market (PDF)"),
    @(19, "Melbourne", "Yarra Trams - No. 11", "7:55am - 8:10am  11/2/2021", "Case used tram no. 11  Start: D16- Harbour Esplanade/ Collins Street  Finish: William Street/ Collins Street #3"),
    @(20, "Melbourne", "Yarra Trams - No. 58", "8:10am - 8:25am  11/2/2021", "Case used tram no. 58  Start: Bourke Street/ William Street #5  Finish: Queen Victoria Market/ Peel Street #9"),
    @(21, "Melbourne", "Yarra Trams - No. 58", "9:40am - 9:55am  11/2/2021", "Case used tram no. 58  Start: Queen Victoria Market/Peel Street #9  Finish: Bourke Street/ William Street #5")
)

foreach ($row in $newRows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}
